$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dadosDeAcesso")

$ws.Range("A7").Value = "ID_0009"
$ws.Range("B7").Value = "André Automatizador"
$ws.Range("C7").Value = "sem email"
$ws.Range("D7").Value = "automacaoteste"

[void]$ws.Activate()
[void]$ws.Range("D7").Select()
